# Updated cryptos list: refresh Price (column D) and Volume(1h) (column E) values,
# and fix row ordering for a few coins that swapped rank position
# (EthereumClassic/Dai rows 26-27, Toncoin/Cosmos rows 28-29, WEMIXToken/Monero rows 42-43).
# NumberFormat is forced to Text ("@") before writing Price values that look like plain
# decimal numbers, so Excel keeps them as text strings instead of auto-converting them to
# floating point numbers (values that already contain more than one "." stay text on their own).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.368.84"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "2.521.03"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.10"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.11"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.527"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.562"
$ws.Range("E9").Value = "  +4.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.31"
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.81"
$ws.Range("E11").Value = "  +6.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0819"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.25"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "2.920.63"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "2.524.05"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.855"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "48.259.30"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.37"
$ws.Range("E19").Value = "  +3.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.65"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.72"
$ws.Range("E22").Value = "  +1.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.85"
$ws.Range("E23").Value = "  +3.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "273.24"
$ws.Range("E24").Value = "  +9.60%  "
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.19"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.31"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.17"
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.145"
$ws.Range("E30").Value = "  +4.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.23"
$ws.Range("E31").Value = "  -1.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.77"
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.97"
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.42"
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.98"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.71"
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.30"
$ws.Range("E41").Value = "  +3.77%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.19"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "118.76"
$ws.Range("E43").Value = "  -2.74%  "
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("D45").Value = "2.001.30"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.12"
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("E47").Value = "  +5.52%  "
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.09"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.34"
$ws.Range("E51").Value = "  +3.61%  "
